# Sign Up page tests were added
# - emails are shortened (first_member@mail.com -> first@mail.com, etc.)
# - the email column (D) becomes a mailto: hyperlink for each user
# - selection moves to D7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the email values for the three data rows.
$ws.Range("D2").Value = "first@mail.com"
$ws.Range("D3").Value = "second@mail.com"
$ws.Range("D4").Value = "third@mail.com"

# Turn each email cell into a mailto hyperlink (this also applies the
# built-in "Hyperlink" cell style, matching the new style added to
# styles.xml).
$null = $ws.Hyperlinks.Add($ws.Range("D2"), "mailto:first@mail.com")
$null = $ws.Hyperlinks.Add($ws.Range("D3"), "mailto:second@mail.com")
$null = $ws.Hyperlinks.Add($ws.Range("D4"), "mailto:third@mail.com")

# Match the selected cell recorded in the saved workbook.
$null = $ws.Range("D7").Select()
